# Commit: "Mon, May 25, 2020  6:05:11 PM"
#
# The underlying OOXML diff has one semantically-reachable change through the
# PowerPoint object model: the table on slide 5 (ppt/slides/slide5.xml) gets
# its <a:tableStyleId> switched from the custom "Table_0" style
# ({BD82FF4E-F017-4849-866B-73DEA91539C4}) to the built-in table style
# {0283E190-1FF2-43F7-92B9-CD9EB65A392C}, while bandRow/firstRow/noFill stay
# untouched.
#
# (The rest of the diff just moves the existing "Integral" / "Office Theme"
# theme parts' raw XML between ppt/theme/theme1.xml and ppt/theme/theme2.xml
# with no relationship/ID changes at all -- there is no Theme/Design COM
# member in this object model that can rewrite clrScheme/fontScheme/fmtScheme
# or the theme's name, so that half of the commit has no COM-reachable
# equivalent and is intentionally left alone here.)

$p = $ppt.ActivePresentation

# slide5.xml is the 5th slide in the deck; it holds a single table inside a
# graphicFrame shape.
$slide = $p.Slides.Item(5)

foreach ($shape in $slide.Shapes) {
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{0283E190-1FF2-43F7-92B9-CD9EB65A392C}")
    }
}
